$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "Unveiling the Complexity of Cosmic Phenomena" "Exploring the Marvelous World of Chemistry: A Journey into the Realm of Elements and Compounds"

Replace-Text " Neil deGrasse Tyson" " Eleanor Stanton"

Replace-Text "NeilTyson@SpaceInstitute" "estanton@edu"

Replace-Text "Embarking on an enchanting voyage to unravel the complexities of cosmic phenomena, we delve into the profound mysteries that enchant our universe" "Chemistry, an intriguing and impactful science, unveils the hidden intricacies of matter and its diverse interactions"

Replace-Text " From the grand tapestry of galaxies that stretch across unfathomable distances, to the enigmatic fabric of space and time warping around celestial bodies, our quest for comprehension leads us down a path of captivating discoveries" " From the vast universe to the microscopic realm within our bodies, chemistry plays a pivotal role in shaping our world"

Replace-Text " Through meticulous observation, tireless calculations, and imaginative leap, humanity continues to unlock the secrets of the cosmos, redefining our perception of existence itself" " As we delve into the fascinating tapestry of chemistry, we embark on an exhilarating exploration of the elements that constitute everything around us and the myriad compounds formed through their intricate combinations"

Replace-Text "As our telescopes peer deeper into the vast expanse, we encounter distant worlds that ignite our curiosity and contemplation" "In this realm of substances, we uncover the fundamental principles governing chemical reactions, witnessing the wondrous transformations of matter into new entities with unique properties"

Replace-Text " The interplay of cosmic forces, the birth and death of stars, and the symphony of interactions between celestial bodies fuel our inquiry into the fundamental principles that govern the universe's evolution" " The symphony of chemistry encompasses myriad concepts, from atomic structures and bonding arrangements to energy transfer and reaction dynamics"

Replace-Text " Each cosmic event, each intricate celestial dance, holds clues to unraveling mysteries that have captivated humankind for millennia, painting a breathtaking tapestry of cosmic beauty and awe" " Each element, with its distinctive characteristics, contributes to the intricate dance of chemical interactions, orchestrating the formation of countless compounds with diverse applications in fields ranging from medicine to materials science"

Replace-Text "Yet, the complexities of the universe extend beyond the reaches of our tangible world, delving into realms that transcend our current understanding" "As we unravel the enigmas of chemistry, we gain invaluable insights into the natural world, unveiling the intricate mechanisms underlying life itself"

Replace-Text " From the mysteries of dark matter and energy, whose enigmatic nature eludes our grasp, to the theoretical concept of multiple universes and the complexities of multi-dimensional space, our exploration into the cosmos opens doors to realms that challenge our conceptual boundaries" " From the intricate workings of photosynthesis, the process by which plants convert sunlight into energy, to the intricate pathways of cellular respiration, the fundamental energy-generating process within living organisms, chemistry unveils the symphony of life at its most fundamental level"

Replace-Text " In this journey of seeking cosmic knowledge, we embark on an intellectual adventure where wonder and enigma intersect, inviting us to the depths of the universe's profound secrets" " Its principles permeate every aspect of our existence, shaping the materials we use, the medicines that heal us, and the intricate complexity of the living world"

Replace-Text "Our journey of exploration into cosmic phenomena unveils a harmonious tapestry of elegance, mystery, and boundless beauty" "This essay embarks on an enthralling exploration of chemistry, venturing into the captivating realm of elements, compounds, and their captivating interactions"

Replace-Text " From the observable grandeur of distant galaxies and the intricacies of gravitational dance to the enigma of unseen forces and the complexities of spacetime, the universe constantly presents us with profound mysteries that test the limits of our comprehension" " From the fundamental "

Replace-Text " As we continue to probe the cosmos with unrelenting curiosity, we are " "principles governing chemical reactions to the intricacies of life itself, chemistry weaves the tapestry of our world, influencing countless aspects of our existence"

Replace-Text "reminded that the path to cosmic understanding is a never-ending adventure, filled with captivating discoveries and awe-inspiring revelations that paint a mesmerizing portrait of the intricate universe we inhabit" ""

Replace-Text " Embracing the enigma of the cosmos, we venture forth with an insatiable hunger for knowledge, eager to unravel the secrets that the universe holds" " Through the study of chemistry, we gain a profound understanding of the natural world, unlocking the secrets of matter and its remarkable transformations, revealing the symphony of life at its most fundamental level"
